# Complete wikidata links in 1701 places
# Strip stray research/debug comment text that had been accidentally
# prepended to several "comment" (column E) cells, and fix the one
# row (26) whose wikidata id/comment were swapped with leftover notes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E5").Value = 'Jen-houo, hoje: Renhe, 仁和县, Historical county name. coordinate: 30.448897N, 120.307504E'
$ws.Range("E6").Value = 'K''iu-tcheou, hoje:Quzhou, 衢州, , in the Chinese translation it is recognized as “遂州”, which is wrong, both phonetically and geographically. In Dehergne(1957), it is noted as "衢州".'
$ws.Range("E10").Value = 'Kia-chan, hoje: Jiashan, 嘉善, Kaosham'
$ws.Range("D26").Value = 'No wikidata'
$ws.Range("E26").Value = ""
$ws.Range("E34").Value = 'Zhangzhou, hoje:Zhangzhou, 漳州, Changchow (Lungki)'
$ws.Range("E35").Value = 'hoje: Houban, 后坂, (@geonames:1977135)Au-poa,Heupuen'
$ws.Range("E65").Value = 'Ki-long, Jilong, hoje:Keelung, 基隆, '
$ws.Range("E67").Value = 'hoje：大包里'
$ws.Range("E70").Value = 'hoje:Henan, 河南, '
$ws.Range("E73").Value = 'Koei-té,Kueite, hoje: Guide, 归德, '
$ws.Range("E74").Value = 'Huguang?, Hou-Quang, hoje:Huguang,湖广, '
$ws.Range("E81").Value = 'do Hupei?, K''i-tcheou,Chichou, hoje: Qizhou, 蕲州, '
$ws.Range("E94").Value = ""
$ws.Range("E112").Value = 'Song-kiang,Sungching, hoje: Songjiang, 松江, '
$ws.Range("E119").Value = 'Chang-hai, hoje: Shanghai, 上海, '
$ws.Range("E121").Value = 'Tsi-pao,Chipao, hoje: Qibao, 七宝, '
$ws.Range("E138").Value = 'Lapa,Wantchai, hoje: Wanzai, 湾仔, Lappa (Wantchai)'
$ws.Range("E143").Value = 'K''iong-tcheou,Chiungchou, hoje: Qiongzhou, 琼州, Kiungchow (île de Hainan)'
$ws.Range("E151").Value = 'Cinçun, hoje: Jingcun, 靖村, coordinate:24.840448198893206N, 113.54394322209676E. In the original book (Dehergne, 1973), it is written as "Tsintsun" (without g), which is wrongly spelled. In Dehergne(1973), it is written as Tsingtsun, with the Chinese name 靖村 noted. It is "à une heure et demie de marche au nord-ouest de Siuchow." In the Chinese traslation, it is recognized as "青村", which is wrong.'
$ws.Range("E152").Value = 'Vançun, hoje: Huangcun, 黄村, coordinate: 24.900778941203768N, 113.60398607001083E Dehergne(1957): "à une heure de chemin au nord de Shiuchow".'
$ws.Range("E153").Value = 'Hwanghsiaping, hoje: Vankaxen, 黄下坪？, In the Chinese translation, it is recognized as "黄下坪". Dehergne(1957):"Hwanghsiaping (" Vankaxen ") à 15 li à l''ouest de Shiuchow." But it cannot be found in the map.'
$ws.Range("E154").Value = 'Mochi, hoje: Madigang, 麻地岗, In the Chinese translation, it is recognized as “末岗”, but in Dehergne(1957), it is noted as “麻地岗”. coordinate: 24.68899887552694N, 113.57599418283718E'
$ws.Range("E155").Value = 'Yang-hiang, hoje: 杨姓村'
$ws.Range("E157").Value = 'Koei-Tcheou, hoje: Guizhou, 贵州, '
$ws.Range("E159").Value = 'hoje: 北京, Peking (Shuntien) (1598), 1601'
$ws.Range("E164").Value = 'Tcheng-ting, hoje: Zhengding, 正定, '
$ws.Range("E188").Value = ""
$ws.Range("E197").Value = 'Dehergne(1957)did not give the Chinese name of Peichingtien.In the Chinese translation, it is recognized as “北辛店村”，but the pronunciation dose not match. In addition, the “北辛店村”of Shandong is far from Dongchang.'
$ws.Range("E203").Value = 'Kiao-tou (ou Tungyiianfang)'
$ws.Range("E209").Value = ""
$ws.Range("E213").Value = 'Mien tcheou, hoje: Mianzhou, 绵州, hoje:  '
